$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("RETAIL_DATA")

# ---------------------------------------------------------------
# Sheet1 (RUNMANAGER): re-label headers, flip C2 yes -> no
# ---------------------------------------------------------------
$ws1.Cells.Item(1,1).Value = "Test Case Name"
$ws1.Cells.Item(1,2).Value = "Test Description"
$ws1.Cells.Item(1,3).Value = "Execute"
$ws1.Cells.Item(1,4).Value = "Priority"
$ws1.Cells.Item(1,5).Value = "Count"

$ws1.Cells.Item(2,1).Value = "LoginLogoutTest"
$ws1.Cells.Item(2,2).Value = "To check whether the user can successfully login and logout"
$ws1.Cells.Item(2,3).Value = "no"

$ws1.Cells.Item(3,1).Value = "newTest"
$ws1.Cells.Item(3,2).Value = "To check this test runs"
$ws1.Cells.Item(3,3).Value = "yes"

[void]$ws1.Range("B10").Select()

# ---------------------------------------------------------------
# Sheet2 (RETAIL_DATA): new blank label row, reshuffled data,
# trailing empty rows removed, new rId hyperlinks, extra column.
# ---------------------------------------------------------------
[void]$ws2.Hyperlinks.Delete()
[void]$ws2.Range("A5:E7").EntireRow.Delete()
[void]$ws2.Rows.Item(2).Insert()
[void]$ws2.Range("E5").Clear()

$ws2.Cells.Item(1,1).Value = "Test Case Name"
$ws2.Cells.Item(1,2).Value = "Execute"
$ws2.Cells.Item(1,3).Value = "Browser"
$ws2.Cells.Item(1,4).Value = "data1"
$ws2.Cells.Item(1,5).Value = "data2"

$ws2.Cells.Item(2,1).Value = "'"
$ws2.Cells.Item(2,2).Value = "'"
$ws2.Cells.Item(2,3).Value = "'"
$ws2.Cells.Item(2,4).Value = "User name"
$ws2.Cells.Item(2,5).Value = "Login Password"

$ws2.Cells.Item(3,1).Value = "LoginLogoutTest"
$ws2.Cells.Item(3,2).Value = "yes"
$ws2.Cells.Item(3,3).Value = "chrome"
$ws2.Cells.Item(3,4).Value = "spcbtest"
$ws2.Cells.Item(3,5).Value = "Asdf@123"

$ws2.Cells.Item(4,1).Value = "newTest"
$ws2.Cells.Item(4,2).Value = "yes"
$ws2.Cells.Item(4,3).Value = "chrome"
$ws2.Cells.Item(4,4).Value = "spcb"
$ws2.Cells.Item(4,5).Value = "Asdf@123"

[void]$ws2.Hyperlinks.Add($ws2.Range("E3"), "mailto:Asdf@123")
$ws2.Range("E3").Style = "Hyperlink"
[void]$ws2.Hyperlinks.Add($ws2.Range("E4"), "mailto:Asdf@123")
$ws2.Range("E4").Style = "Hyperlink"

$ws2.Columns.Item(6).ColumnWidth = 15.1667

[void]$ws2.Range("F4").Select()

# ---------------------------------------------------------------
# Make RETAIL_DATA the active tab/sheet (workbook activeTab + tabSelected)
# ---------------------------------------------------------------
[void]$ws2.Activate()

Write-Output "done"
